$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the "core" complex-chart question codes (and the charting-date
#    question code) so their IDs are constant instead of test-only stand-ins.
#    These live in the shared string table and are referenced from the
#    "Core" and "Test Chart" sheets.
# ---------------------------------------------------------------------------
$core = $wb.Worksheets.Item("Core")
$core.Range("A2").Value = "ComplexChartInstanceName"
$core.Range("A3").Value = "ComplexChartDate"
$core.Range("A4").Value = "ComplexChartType"
$core.Range("A5").Value = "ComplexChartSubtype"

$testChart = $wb.Worksheets.Item("Test Chart")
$testChart.Range("A2").Value = "PatientChartingDate"

# ---------------------------------------------------------------------------
# 2. The renamed "code" cells (column A) pick up the alignment formatting
#    that was previously only used by the "type" cells (column B) in the
#    same rows, while those "type" cells (and R5) fall back to the plain,
#    un-aligned formatting. Replicate this by copying formats between the
#    relevant cells (format-painter equivalent), which reassigns the
#    underlying cell style records the same way Excel's UI would.
# ---------------------------------------------------------------------------

# Core!A2:A5 take on the alignment formatting currently used by Core!B2.
$core.Range("B2").Copy()
$core.Range("A2:A5").PasteSpecial(-4122)  # xlPasteFormats

# Core!B2 / Core!B3 / Core!R5 revert to the plain formatting used by Core!C2.
$core.Range("C2").Copy()
$core.Range("B2").PasteSpecial(-4122)     # xlPasteFormats
$core.Range("C2").Copy()
$core.Range("B3").PasteSpecial(-4122)     # xlPasteFormats
$core.Range("C2").Copy()
$core.Range("R5").PasteSpecial(-4122)     # xlPasteFormats

# Test Chart!A2 takes on that same alignment formatting.
$core.Range("B2").Copy()
$testChart.Range("A2").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = $false
